# Update the financial package workbook:
#  - Balance Sheet: add "Prepaid expenses" + "Total Current Assets" rows (a
#    current-assets subtotal) before "Total Assets", and add "Accounts
#    Payable", "Deferred revenue" and "Total Equity" rows after "Total
#    Assets".
#  - Income Statement: add "Gross Profit", "Salaries" and "Total Operating
#    expenses" rows before "Net Income".
#
# The Amount/Variance/% Variance columns hold numeric-looking text (e.g.
# "$1,396,603.84", "-2.75%") that must stay literal text, matching the rest
# of the workbook, rather than being auto-parsed into numbers. Formatting
# the range as Text before assigning the values keeps them literal; the
# style is then reset to Normal so the cells end up unstyled, just like
# every other cell in these sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Balance Sheet
# ---------------------------------------------------------------------------
$bs = $wb.Worksheets.Item("Balance Sheet")

# Row 4 currently holds "Total Assets". Insert two new rows above it for
# "Prepaid expenses" (new row 4) and "Total Current Assets" (new row 5);
# "Total Assets" is pushed down to row 6.
$bs.Rows.Item(4).Insert()
$bs.Rows.Item(4).Insert()

$bs.Range("B4:E5").NumberFormat = "@"

$bs.Range("A4").Value = "Prepaid expenses"
$bs.Range("B4").Value = "$1,396,603.84"
$bs.Range("C4").Value = "$1,436,064.00"
$bs.Range("D4").Value = "($39,460.16)"
$bs.Range("E4").Value = "-2.75%"

$bs.Range("A5").Value = "Total Current Assets"
$bs.Range("B5").Value = "$28,086,600.03"
$bs.Range("C5").Value = "$29,493,770.49"
$bs.Range("D5").Value = "($1,407,170.46)"
$bs.Range("E5").Value = "-4.77%"

$bs.Range("A4:E5").Style = "Normal"

# "Total Assets" is now row 6. Append the liabilities/equity rows after it
# (rows 7-9 are currently blank, so no insert is needed).
$bs.Range("B7:E9").NumberFormat = "@"

$bs.Range("A7").Value = "Accounts Payable"
$bs.Range("B7").Value = "$588,450.01"
$bs.Range("C7").Value = "$333,741.90"
$bs.Range("D7").Value = "$254,708.11"
$bs.Range("E7").Value = "76.32%"

$bs.Range("A8").Value = "Deferred revenue"
$bs.Range("B8").Value = "$2,343,160.23"
$bs.Range("C8").Value = "$2,281,494.61"
$bs.Range("D8").Value = "$61,665.62"
$bs.Range("E8").Value = "2.70%"

$bs.Range("A9").Value = "Total Equity"
$bs.Range("B9").Value = "$24,403,574.05"
$bs.Range("C9").Value = "$26,254,039.93"
$bs.Range("D9").Value = "($1,850,465.88)"
$bs.Range("E9").Value = "-7.05%"

$bs.Range("A7:E9").Style = "Normal"

# ---------------------------------------------------------------------------
# Income Statement
# ---------------------------------------------------------------------------
$incSt = $wb.Worksheets.Item("Income Statement")

# Row 5 currently holds "Net Income". Insert three new rows above it for
# "Gross Profit", "Salaries" and "Total Operating expenses"; "Net Income"
# is pushed down to row 8.
$incSt.Rows.Item(5).Insert()
$incSt.Rows.Item(5).Insert()
$incSt.Rows.Item(5).Insert()

$incSt.Range("B5:E7").NumberFormat = "@"

$incSt.Range("A5").Value = "Gross Profit"
$incSt.Range("B5").Value = "$1,223,682.10"
$incSt.Range("C5").Value = "$1,145,606.48"
$incSt.Range("D5").Value = "$78,075.62"
$incSt.Range("E5").Value = "6.82%"

$incSt.Range("A6").Value = "Salaries"
$incSt.Range("B6").Value = "$1,380,987.81"
$incSt.Range("C6").Value = "$1,326,380.04"
$incSt.Range("D6").Value = "$54,607.77"
$incSt.Range("E6").Value = "4.12%"

$incSt.Range("A7").Value = "Total Operating expenses"
$incSt.Range("B7").Value = "$3,224,657.83"
$incSt.Range("C7").Value = "$2,919,932.52"
$incSt.Range("D7").Value = "$304,725.31"
$incSt.Range("E7").Value = "10.44%"

$incSt.Range("A5:E7").Style = "Normal"
